$wb = $excel.ActiveWorkbook

# ALC!row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2207.6562
$ws.Range("I15").Value = 2207.6562
$ws.Range("K15").Value = 6622.9686
$ws.Range("M15").Value = -6453.9686

# ALC!row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 166667090
$ws.Range("I18").Value = 503.2
$ws.Range("K18").Value = 503.2
$ws.Range("M18").Value = -219.2

# ALC!row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 189.55556
$ws.Range("I33").Value = 212.57143
$ws.Range("K33").Value = 212.57143
$ws.Range("M33").Value = 16.42857000000001

# ALC!row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 161.71428
$ws.Range("J41").Value = 168.6
$ws.Range("L41").Value = 168.6
$ws.Range("N41").Value = -1048.6

# ALC!row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 7527.6665
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 7527.6665
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 22582.9995
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -23122.9995

# ALC!row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 7527.6665
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 7527.6665
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 22582.9995
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -24454.9995

# ALC!row 97
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 2299.2
$ws.Range("J97").Value = 2299.2
$ws.Range("L97").Value = 6897.599999999999
$ws.Range("N97").Value = -7889.599999999999

# ALC!row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4728.7617
$ws.Range("I98").Value = 4296.4443
$ws.Range("K98").Value = 4296.4443
$ws.Range("M98").Value = -2798.4443

# ALC!row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 4763300
$ws.Range("I103").Value = 2030.5555
$ws.Range("K103").Value = 6091.666499999999
$ws.Range("M103").Value = -5505.666499999999

# ALC!row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4327.4287
$ws.Range("I116").Value = 4325.273
$ws.Range("K116").Value = 4325.273
$ws.Range("M116").Value = -883.2730000000001

# ALC!row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 4728.7617
$ws.Range("I122").Value = 4296.4443
$ws.Range("K122").Value = 12889.3329
$ws.Range("M122").Value = -10439.3329

# ALC!row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 564012.3
$ws.Range("I131").Value = 722765.6
$ws.Range("K131").Value = 2168296.8
$ws.Range("M131").Value = -2163256.8

# ALC!row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4476.1143
$ws.Range("I132").Value = 4578.353
$ws.Range("K132").Value = 13735.059
$ws.Range("M132").Value = -11205.059

# ALC!row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1356579.4
$ws.Range("I137").Value = 1924872.6
$ws.Range("K137").Value = 5774617.800000001
$ws.Range("M137").Value = -5772067.800000001

# ARM!row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3478.6553
$ws.Range("I32").Value = 3495.8215
$ws.Range("K32").Value = 3495.8215
$ws.Range("M32").Value = -3208.8215

# ARM!row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2582.4211
$ws.Range("I61").Value = 1230.1
$ws.Range("K61").Value = 1230.1
$ws.Range("M61").Value = -1018.1

# ARM!row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2292.7585
$ws.Range("I122").Value = 2194.5
$ws.Range("J122").Value = 3144.3333
$ws.Range("K122").Value = 6583.5
$ws.Range("L122").Value = 9432.999899999999
$ws.Range("M122").Value = -4133.5
$ws.Range("N122").Value = -14332.9999

# ARM!row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2248.3684
$ws.Range("I132").Value = 1376.9
$ws.Range("K132").Value = 4130.700000000001
$ws.Range("M132").Value = -1600.700000000001

# ARM!row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2582.4211
$ws.Range("I136").Value = 1230.1
$ws.Range("K136").Value = 3690.3
$ws.Range("M136").Value = -1140.3

# BSM!row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 11305680
$ws.Range("I105").Value = 715358.5
$ws.Range("J105").Value = 27779514
$ws.Range("K105").Value = 715358.5
$ws.Range("L105").Value = 27779514
$ws.Range("M105").Value = -713611.5
$ws.Range("N105").Value = -27783008

# CRP!row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4038.8867
$ws.Range("I31").Value = 3023.3447
$ws.Range("K31").Value = 3023.3447
$ws.Range("M31").Value = -2728.3447

# CRP!row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4038.8867
$ws.Range("I34").Value = 3023.3447
$ws.Range("K34").Value = 3023.3447
$ws.Range("M34").Value = -2821.3447

# CRP!row 54
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

# CRP!row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2498.3333
$ws.Range("I58").Value = 1907.2
$ws.Range("K58").Value = 1907.2
$ws.Range("M58").Value = -1704.2

# CRP!row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2045.6154
$ws.Range("I122").Value = 2124.6
$ws.Range("K122").Value = 6373.799999999999
$ws.Range("M122").Value = -3923.799999999999

# CRP!row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 12825570
$ws.Range("I132").Value = 4162.643
$ws.Range("K132").Value = 12487.929
$ws.Range("M132").Value = -9957.929

# CRP!row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4300.9
$ws.Range("I134").Value = 4664
$ws.Range("K134").Value = 13992
$ws.Range("M134").Value = -11457

# CRP!row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2498.3333
$ws.Range("I136").Value = 1907.2
$ws.Range("K136").Value = 5721.6
$ws.Range("M136").Value = -3171.6

# CUL!row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 506.875
$ws.Range("J92").Value = 537.8570999999999
$ws.Range("L92").Value = 1613.5713
$ws.Range("N92").Value = -4109.5713

# CUL!row 124
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 5287.5
$ws.Range("I124").Value = 575
$ws.Range("J124").Value = 10000
$ws.Range("K124").Value = 1725
$ws.Range("L124").Value = 30000
$ws.Range("M124").Value = 3185
$ws.Range("N124").Value = -39820

# CUL!row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 2768.4285
$ws.Range("J134").Value = 4996
$ws.Range("L134").Value = 14988
$ws.Range("N134").Value = -25128

# CUL!row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 5421.5
$ws.Range("J136").Value = 13500
$ws.Range("L136").Value = 40500
$ws.Range("N136").Value = -50700

# CUL!row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 55557920
$ws.Range("I137").Value = 83335050
$ws.Range("J137").Value = 3665.3333
$ws.Range("K137").Value = 250005150
$ws.Range("L137").Value = 10995.9999
$ws.Range("M137").Value = -250000050
$ws.Range("N137").Value = -21195.9999

# GSM!row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 58827190
$ws.Range("I80").Value = 100003100
$ws.Range("J80").Value = 4451.4287
$ws.Range("K80").Value = 100003100
$ws.Range("L80").Value = 4451.4287
$ws.Range("M80").Value = -100002102
$ws.Range("N80").Value = -6447.4287

# GSM!row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 58827190
$ws.Range("I83").Value = 100003100
$ws.Range("J83").Value = 4451.4287
$ws.Range("K83").Value = 500015500
$ws.Range("L83").Value = 22257.1435
$ws.Range("M83").Value = -500010508
$ws.Range("N83").Value = -32241.1435

# GSM!row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1523.579
$ws.Range("I102").Value = 1089.1111
$ws.Range("J102").Value = 1914.6
$ws.Range("K102").Value = 1089.1111
$ws.Range("L102").Value = 1914.6
$ws.Range("M102").Value = 532.8888999999999
$ws.Range("N102").Value = -5158.6

# GSM!row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2564.375
$ws.Range("I113").Value = 2497.25
$ws.Range("J113").Value = 2631.5
$ws.Range("K113").Value = 2497.25
$ws.Range("L113").Value = 2631.5
$ws.Range("M113").Value = -327.25
$ws.Range("N113").Value = -6971.5

# LTW!row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1685.0769
$ws.Range("I16").Value = 538.1667
$ws.Range("J16").Value = 2668.1428
$ws.Range("K16").Value = 538.1667
$ws.Range("L16").Value = 2668.1428
$ws.Range("M16").Value = -368.1667
$ws.Range("N16").Value = -3008.1428

# LTW!row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 25001
$ws.Range("I40").Value = 25001
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 25001
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -24865
$ws.Range("N40").ClearContents()

# LTW!row 42
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 32014
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 32014
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 32014
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -33140

# LTW!row 49
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H49").Value = 32014
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 32014
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 32014
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -32308

# LTW!row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3299
$ws.Range("I68").Value = 2949
$ws.Range("J68").Value = 3999
$ws.Range("K68").Value = 2949
$ws.Range("L68").Value = 3999
$ws.Range("M68").Value = -2200
$ws.Range("N68").Value = -5497

# LTW!row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3299
$ws.Range("I71").Value = 2949
$ws.Range("J71").Value = 3999
$ws.Range("K71").Value = 14745
$ws.Range("L71").Value = 19995
$ws.Range("M71").Value = -11001
$ws.Range("N71").Value = -27483

# LTW!row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 410.85715
$ws.Range("I93").Value = 373.72726
$ws.Range("J93").Value = 547
$ws.Range("K93").Value = 373.72726
$ws.Range("L93").Value = 547
$ws.Range("M93").Value = 874.27274
$ws.Range("N93").Value = -3043

# WVR!row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3650
$ws.Range("I126").Value = 3650
$ws.Range("K126").Value = 10950
$ws.Range("M126").Value = -8480

# WVR!row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2109.0588
$ws.Range("I132").Value = 1844
$ws.Range("J132").Value = 3346
$ws.Range("K132").Value = 5532
$ws.Range("L132").Value = 10038
$ws.Range("M132").Value = -3002
$ws.Range("N132").Value = -15098
